# Tidy up of switch to holder class
#
# Adds a new "holder class" results row (row 6) to the "G4 (Iteration 1)" and
# "Pi4 (TreeMap)" sheets, extends the shared % Improvement formula down to
# cover it, updates the "Pi4 (TreeMap)" bar chart's series ranges to include
# the new point, and leaves the cursor where the author left it on each
# sheet.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "G4 (Iteration 1)": new row 6 = "4 (class)" holder-class result
# ------------------------------------------------------------------
$wsG4 = $wb.Worksheets.Item("G4 (Iteration 1)")

$wsG4.Range("A6").Value = "4 (class)"
$wsG4.Range("B6").Value = "FastDTW_1NN"
$wsG4.Range("C6").Value = "ChlorineConcentration"
$wsG4.Range("D6").Value = 0
$wsG4.Range("E6").Value = 2828352.2450000001
$wsG4.Range("E6").Style = "Normal"
$wsG4.Range("F6").Value = 10860872620
$wsG4.Range("F6").Style = "Normal"
$wsG4.Range("G6").Value = 1568099091
$wsG4.Range("G6").Style = "Normal"
$wsG4.Range("H6").Formula = "=100*(1-(E6/`$E`$2))"
$wsG4.Range("H6").Style = "Normal"

# Cursor left selecting the new row after entering the data
$wsG4.Activate()
$wsG4.Range("A7:H7").Select()

# ------------------------------------------------------------------
# Sheet "Pi4 (TreeMap)": new row 6 = "Vars (4)" holder-class result
# ------------------------------------------------------------------
$wsTreeMap = $wb.Worksheets.Item("Pi4 (TreeMap)")

$wsTreeMap.Range("A6").Value = "Vars (4)"
$wsTreeMap.Range("B6").Value = "FastDTW_1NN"
$wsTreeMap.Range("C6").Value = "ChlorineConcentration"
$wsTreeMap.Range("D6").Value = -1
$wsTreeMap.Range("E6").Value = 26822213.548999999
$wsTreeMap.Range("E6").NumberFormat = "0.00"
$wsTreeMap.Range("F6").Value = 102997300027
$wsTreeMap.Range("F6").NumberFormat = "0.0"
$wsTreeMap.Range("G6").Value = 12038206017
$wsTreeMap.Range("G6").NumberFormat = "0.0"
$wsTreeMap.Range("H6").Formula = "=100*(1-(E6/`$E`$2))"
$wsTreeMap.Range("H6").Style = "Normal"

# Update the bar chart on this sheet so its series covers the new row
$co = $wsTreeMap.ChartObjects().Item(1)
$chart = $co.Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(,'Pi4 (TreeMap)'!`$A`$2:`$A`$6,'Pi4 (TreeMap)'!`$E`$2:`$E`$6,1)"

# Cursor left on this (the active) sheet after the edit
$wsTreeMap.Activate()
$wsTreeMap.Range("H26").Select()
